$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Macro_taxonomy")
$ws2.Range("C18").Value = "EWV/LN"
$ws2.Range("C33").Value = "EWV/LN"

$ws7 = $wb.Worksheets.Item("Height")
$ws7.Range("A13").Value = "EWV/LN"
$null = $ws7.Rows("14").Insert()
$ws7.Range("A14").Value = "MATO"
$ws7.Range("B14").Value = "H:1"
$ws7.Range("C14").Value = 1

$ws2.Activate()
$null = $ws2.Range("B24").Select()
